$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price cells that are purely numeric-looking strings stay as text
# (matches original file where all Price column values are stored as text)
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D19", "D20", "D21", "D22", "D23", "D24", "D28", "D29", "D30", "D31", "D33", "D34", "D38", "D39", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "65.685.97"
$ws.Range("E2").Value = "  -1.69%  "
$ws.Range("D3").Value = "3.409.37"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "596.74"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").Value = "142.58"
$ws.Range("E6").Value = "  -3.69%  "
$ws.Range("B7").Value = "LidoStakedEther"
$ws.Range("C7").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D7").Value = "3.407.91"
$ws.Range("E7").Value = "  -1.54%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "0.470"
$ws.Range("E9").Value = "  -2.82%  "
$ws.Range("D10").Value = "0.135"
$ws.Range("E10").Value = "  -5.18%  "
$ws.Range("D11").Value = "7.94"
$ws.Range("E11").Value = "  +5.81%  "
$ws.Range("D12").Value = "0.406"
$ws.Range("E12").Value = "  -4.29%  "
$ws.Range("D13").Value = "3.984.89"
$ws.Range("D14").Value = "0.0000201"
$ws.Range("E14").Value = "  -6.39%  "
$ws.Range("D15").Value = "29.88"
$ws.Range("E15").Value = "  -5.98%  "
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "3.403.47"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").Value = "65.625.88"
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("D19").Value = "10.45"
$ws.Range("E19").Value = "  +4.35%  "
$ws.Range("D20").Value = "6.12"
$ws.Range("E20").Value = "  -5.23%  "
$ws.Range("D21").Value = "14.84"
$ws.Range("E21").Value = "  -3.61%  "
$ws.Range("D22").Value = "417.63"
$ws.Range("E22").Value = "  -5.12%  "
$ws.Range("D23").Value = "0.582"
$ws.Range("E23").Value = "  -4.64%  "
$ws.Range("D24").Value = "77.54"
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").Value = "3.545.42"
$ws.Range("E26").Value = "  -1.45%  "
$ws.Range("E27").Value = "  -8.30%  "
$ws.Range("D28").Value = "9.30"
$ws.Range("E28").Value = "  -5.94%  "
$ws.Range("D29").Value = "7.86"
$ws.Range("E29").Value = "  -6.96%  "
$ws.Range("D30").Value = "2.43"
$ws.Range("E30").Value = "  -2.28%  "
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  -3.94%  "
$ws.Range("D33").Value = "1.48"
$ws.Range("E33").Value = "  -8.29%  "
$ws.Range("D34").Value = "24.52"
$ws.Range("E34").Value = "  -3.63%  "
$ws.Range("D35").Value = "3.403.14"
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  -5.71%  "
$ws.Range("D38").Value = "5.58"
$ws.Range("E38").Value = "  -8.44%  "
$ws.Range("D39").Value = "7.59"
$ws.Range("E39").Value = "  -4.26%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "169.58"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").Value = "0.0864"
$ws.Range("E42").Value = "  -3.17%  "
$ws.Range("D43").Value = "5.09"
$ws.Range("E43").Value = "  -6.09%  "
$ws.Range("D44").Value = "0.872"
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("D45").Value = "1.93"
$ws.Range("E45").Value = "  -10.96%  "
$ws.Range("D46").Value = "45.52"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("D47").Value = "26.78"
$ws.Range("E47").Value = "  -8.11%  "
$ws.Range("D48").Value = "1.18"
$ws.Range("E48").Value = "  -5.41%  "
$ws.Range("D49").Value = "7.11"
$ws.Range("E49").Value = "  -5.01%  "
$ws.Range("D50").Value = "2.31"
$ws.Range("E50").Value = "  -6.56%  "
$ws.Range("D51").Value = "0.927"
$ws.Range("E51").Value = "  -6.03%  "
